# This script applies spell/grammar-correction edits to the "Rpi vs Arduino"
# and "Navio2 vs Pixhawk" comparison tables in the task-list document.
#
# Word's spell/grammar checker, when it autocorrects a misspelled word
# inside a run, splits that run into several runs (one per corrected/kept
# text fragment) and brackets the corrected fragment with <w:proofErr>
# markers (spellStart/spellEnd or gramStart/gramEnd). We reproduce that
# here by rebuilding each affected paragraph's contents via Range.InsertXML
# with the exact desired run/proofErr/bookmark structure. InsertXML always
# replaces the whole paragraph that intersects the target range when the
# payload contains block-level (<w:p>) content, so we always select the
# complete paragraph range before calling it.

$d = $word.ActiveDocument

function Set-ParagraphXml($startPos, $endPos, $pAttrs, $pPrXml, $bodyXml) {
    $range = $d.Range($startPos, $endPos)
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p' + $pAttrs + '>' + $pPrXml + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $range.InsertXML($xml)
}

# Locate each target paragraph by its (old) distinguishing text so the
# script is resilient to minor offset differences.
function Find-ParagraphRange($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return @($p.Range.Start, $p.Range.End)
        }
    }
    return $null
}

# 1) "Wifi capable" -> "Wifi" (wrapped in spellStart/spellEnd) + " capable"
$pos = Find-ParagraphRange("Wifi capable")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00877C24" w:rsidRDefault="00877C24" w:rsidP="00877C24"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' ('<w:proofErr w:type="spellStart"/><w:r><w:t>Wifi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> capable</w:t></w:r>')

# 2) "Diffucult to power using battery pack" -> "Difficult" + " to power using battery pack"
$pos = Find-ParagraphRange("to power using battery pack")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00877C24" w:rsidRDefault="00877C24" w:rsidP="00877C24"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' ('<w:r><w:t>Difficult</w:t></w:r><w:r><w:t xml:space="preserve"> to power using battery pack</w:t></w:r>')

# 3) "Run on linux environment" -> "Run on " + "Linux" + " environment"
$pos = Find-ParagraphRange("Run on linux environment")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00DA6C28" w:rsidRDefault="00163668" w:rsidP="00163668"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">Run on </w:t></w:r><w:r><w:t>Linux</w:t></w:r><w:r><w:t xml:space="preserve"> environment</w:t></w:r>')

# 4) "Easy to setup and intereact with" -> "Easy to setup and " + "interact" + " with"
$pos = Find-ParagraphRange("Easy to setup and intereact with")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00163668"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">Easy to setup and </w:t></w:r><w:r><w:t>interact</w:t></w:r><w:r><w:t xml:space="preserve"> with</w:t></w:r>')

# 5) "Can update firmware wirelessly ( with wifi ) " -> "...with " + "Wi-Fi" + " ) "
$pos = Find-ParagraphRange("Can update firmware wirelessly")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00163668"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">Can update firmware wirelessly ( with </w:t></w:r><w:r><w:t>Wi-Fi</w:t></w:r><w:r><w:t xml:space="preserve"> ) </w:t></w:r>')

# 6) "Very new and not widly seen in industiral applications" -> split + fix spelling;
#    the _GoBack bookmark that used to sit at the end of this paragraph moves to
#    paragraph 9 below (that's where the author's last edit in the source ended up).
$pos = Find-ParagraphRange("Very new and not widly seen")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00163668"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">Very new and not </w:t></w:r><w:r><w:t>widely</w:t></w:r><w:r><w:t xml:space="preserve"> seen in </w:t></w:r><w:r><w:t>industrial</w:t></w:r><w:r><w:t xml:space="preserve"> applications</w:t></w:r>')

# 7) "More consistant in peforming tasks" -> "More " + "consistent" + " in " + "performing" + " tasks"
$pos = Find-ParagraphRange("More consistant in peforming tasks")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00DA6C28"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">More </w:t></w:r><w:r><w:t>consistent</w:t></w:r><w:r><w:t xml:space="preserve"> in </w:t></w:r><w:r><w:t>performing</w:t></w:r><w:r><w:t xml:space="preserve"> tasks</w:t></w:r>')

# 8) "Dicated computing power for flught control firmware" -> "Dictated" + " computing power for " + "flight" + " control firmware"
$pos = Find-ParagraphRange("Dicated computing power for flught control firmware")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00DA6C28"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t>Dictated</w:t></w:r><w:r><w:t xml:space="preserve"> computing power for </w:t></w:r><w:r><w:t>flight</w:t></w:r><w:r><w:t xml:space="preserve"> control firmware</w:t></w:r>')

# 9) "Can be supplemented with external computing power ( like Raspberry Pi )" ->
#    split into runs, move the _GoBack bookmark here, drop the space after "(",
#    and wrap "Pi )" with gramStart/gramEnd.
$pos = Find-ParagraphRange("Can be supplemented with external computing power")
Set-ParagraphXml $pos[0] $pos[1] ' w:rsidR="00163668" w:rsidRDefault="00163668" w:rsidP="00DA6C28"' '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:jc w:val="both"/></w:pPr>' ('<w:r><w:t xml:space="preserve">Can be supplemented with external computing power </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>(like</w:t></w:r><w:r><w:t xml:space="preserve"> Raspberry </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Pi )</w:t></w:r><w:proofErr w:type="gramEnd"/>')

Write-Output "Done"
